$wb = $excel.ActiveWorkbook

# OFF sheet (Week 13 update) - Row 3 ("R")
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B3").Value = 296
$wsOff.Range("C3").Value = 190
$wsOff.Range("D3").Value = 68
$wsOff.Range("E3").Value = 33

# DEF sheet (Week 13 update) - Row 3 ("R")
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B3").Value = 363
$wsDef.Range("C3").Value = 251
$wsDef.Range("D3").Value = 79
$wsDef.Range("F3").Value = 7
